# Update benchmark: 2026-02-10 07:13:12 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI
$ws.Range("F2").Value = "33,33 TL - 33,33 TL"

# Row 3 - HESAPTAN EFT - Şube  (value moved from F to C)
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = ""

# Row 4 - HESAPTAN EFT - ATM  (value moved from F to C)
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = ""

# Row 5 - HESAPTAN EFT - Mobil  (value moved from F to C)
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = ""

# Row 6 - DÜZENLİ EFT
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"

# Row 8 - HESAPTAN HAVALE - Şube  (value moved from F to C, new C value)
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F8").Value = ""

# Row 9 - HESAPTAN HAVALE - ATM  (value moved from F to C, new C value)
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F9").Value = ""

# Row 10 - HESAPTAN HAVALE - Mobil  (value moved from F to C, new C value)
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("F10").Value = ""

# Row 11 - DÜZENLİ HAVALE
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 795 TL | Azami 4.005 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("F14").Value = "2.785,72 TL - 12.380,95 TL"

# Row 15 - ÇEK TAHSİLİ BAŞKA BANKA
$ws.Range("F15").Value = "%0,5 Asgari Tutar: 361,9 TL Azami Tutar: 361,9 TL / 361,9 TL"

# Row 17 - AYNI ŞUBE ÇEK TAHSİLATI
$ws.Range("F17").Value = "%0,5 Asgari Tutar: 427,62 TL Azami Tutar: 427,62 TL"

# Row 20 - ÇEK İADE
$ws.Range("F20").Value = "123,81 TL"

# Row 22 - YP ÇEK TAKASA GÖNDERME
$ws.Range("F22").Value = "%0,5 Asgari Tutar: 427,62 TL Azami Tutar: 427,62 TL / 1.669,52 TL"

# Row 24 - SENET TAHSİLE ALMA
$ws.Range("F24").Value = "600 TL"

# Row 25 - MUAMELESİZ SENET İADESİ
$ws.Range("F25").Value = "600 TL"
